$d = $word.ActiveDocument

# Locate the paragraph that currently holds "Sincerely,"
$p = $d.Paragraphs.Item(21)
$r = $p.Range

# Make room: insert two new (empty) paragraphs right after it.
$r.InsertParagraphAfter()
$r.InsertParagraphAfter()

# --- Paragraph 2 (brand-new): "Let’s check now", highlighted green ---
$p2 = $d.Paragraphs.Item(22)
$p2.Range.Text = "Let’s check now"
$p2.Range.HighlightColorIndex = 4

# --- Paragraph 3 (brand-new): "Sincerely," preceded by a rendered page-break marker ---
$p3 = $d.Paragraphs.Item(23)
$p3xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:lastRenderedPageBreak/><w:t>Sincerely,</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p3.Range.InsertXML($p3xml)

# --- Paragraph 1 (the original "Sincerely," paragraph): repurposed as the testing banner ---
$p1 = $d.Paragraphs.Item(21)
$p1.Range.Text = "Changes made for testing"
$p1.Range.HighlightColorIndex = 7
